$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 196
$ws.Range("F7").Value = 4255
$ws.Range("G7").Value = 85
$ws.Range("G9").Value = 85
$ws.Range("F11").Value = 6235
$ws.Range("G11").Value = 90
$ws.Range("F12").Value = 6235
$ws.Range("G12").Value = 90
$ws.Range("F13").Value = 71
$ws.Range("F15").Value = 2394
$ws.Range("F18").Value = 499
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 9375
$ws.Range("F25").Value = 2503
$ws.Range("F34").Value = 285
$ws.Range("F36").Value = 94
$ws.Range("F41").Value = 106
$ws.Range("F45").Value = 938
$ws.Range("F46").Value = 320
$ws.Range("F48").Value = 33
$ws.Range("F50").Value = 12

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 121

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 703
$ws.Range("F3").Value = 918

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 196
$ws.Range("F3").Value = 703
$ws.Range("F4").Value = 918
$ws.Range("F11").Value = 4255
$ws.Range("G11").Value = 85
$ws.Range("G13").Value = 85
$ws.Range("F16").Value = 6235
$ws.Range("G16").Value = 90
$ws.Range("F17").Value = 71
$ws.Range("F18").Value = 2394
$ws.Range("F20").Value = 499
$ws.Range("F21").Value = 9375
$ws.Range("F26").Value = 2503
$ws.Range("F34").Value = 285
$ws.Range("F36").Value = 94
$ws.Range("F41").Value = 106
$ws.Range("F44").Value = 938
$ws.Range("F45").Value = 320
$ws.Range("F48").Value = 33
$ws.Range("F50").Value = 121
$ws.Range("F51").Value = 121

